# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.909.74"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").Value = "3.453.34"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.84"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.63"
$ws.Range("E6").Value = "  -4.60%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -4.49%  "
$ws.Range("D9").Value = "3.448.61"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  -6.74%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.408"
$ws.Range("E12").Value = "  -4.89%  "
$ws.Range("D13").Value = "4.041.78"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -6.33%  "
$ws.Range("D16").Value = "65.985.45"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("E17").Value = "  -3.78%  "
$ws.Range("D18").Value = "3.447.61"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.89"
$ws.Range("E19").Value = "  -4.93%  "
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.01"
$ws.Range("E21").Value = "  -7.11%  "
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.85"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.533"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").Value = "  -7.41%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.98"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("E32").Value = "  -6.48%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -8.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.01"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.53"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.56"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "29.07"
$ws.Range("E38").Value = "  +11.04%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("D41").Value = "2.746.51"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").Value = "  -10.44%  "
$ws.Range("E43").Value = "  -6.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.32"
$ws.Range("E44").Value = "  -6.55%  "
$ws.Range("E45").Value = "  -4.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.89"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.22"
$ws.Range("E47").Value = "  -7.88%  "
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "304.23"
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.816"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.101"
$ws.Range("E51").Value = "  -3.66%  "
